# Updated cryptos list — apply the price/volume refresh described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Cells whose new "Price" (column D) text looks like a plain decimal number
# (e.g. "64.99") need to be forced to Text first, otherwise the COM layer's
# auto-type-inference would coerce them into a numeric cell (and silently
# drop significant trailing zeros, e.g. "1.00" -> 1). Values that contain
# more than one '.' (thousands-grouped prices like "62.035.58") or other
# non-numeric characters are never coerced, so they don't need this.
# ---------------------------------------------------------------------------
$textCells = @(
    "D5","D6","D13","D14","D19","D20","D21","D22","D24","D25","D27","D28",
    "D31","D32","D33","D34","D36","D37","D38","D40","D42","D43","D44","D45",
    "D47","D48","D49","D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# ---------------------------------------------------------------------------
# Row-by-row value updates: (D = Price, E = Volume(1h))
# ---------------------------------------------------------------------------
$ws.Range("D2").Value  = "62.035.58"
$ws.Range("E2").Value  = "  +2.57%  "

$ws.Range("D3").Value  = "2.423.69"
$ws.Range("E3").Value  = "  +3.83%  "

$ws.Range("E4").Value  = "  +0.09%  "

$ws.Range("D5").Value  = "556.62"
$ws.Range("E5").Value  = "  +2.04%  "

$ws.Range("D6").Value  = "138.40"
$ws.Range("E6").Value  = "  +5.29%  "

$ws.Range("E7").Value  = "  +0.05%  "

$ws.Range("E8").Value  = "  +1.08%  "

$ws.Range("D9").Value  = "2.420.86"
$ws.Range("E9").Value  = "  +3.82%  "

$ws.Range("E10").Value = "  +2.41%  "

$ws.Range("E11").Value = "  +4.57%  "

$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  +3.59%  "

$ws.Range("D14").Value = "25.70"
$ws.Range("E14").Value = "  +8.49%  "

$ws.Range("D15").Value = "2.855.63"
$ws.Range("E15").Value = "  +4.01%  "

$ws.Range("D16").Value = "61.981.67"
$ws.Range("E16").Value = "  +2.58%  "

$ws.Range("E17").Value = "  +4.96%  "

$ws.Range("D18").Value = "2.425.74"
$ws.Range("E18").Value = "  +4.08%  "

$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  +4.69%  "

$ws.Range("D20").Value = "344.66"

$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").Value = "6.81"
$ws.Range("E22").Value = "  +2.44%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "64.99"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").Value = "1.51"
$ws.Range("E27").Value = "  +10.50%  "

$ws.Range("D28").Value = "8.26"
$ws.Range("E28").Value = "  +5.17%  "

$ws.Range("E29").Value = "  +13.11%  "

$ws.Range("D30").Value = "0.0₃0788"
$ws.Range("E30").Value = "  +7.34%  "

$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  +3.88%  "

$ws.Range("D32").Value = "6.32"
$ws.Range("E32").Value = "  +6.12%  "

$ws.Range("D33").Value = "170.84"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").Value = "1.43"
$ws.Range("E34").Value = "  +4.06%  "

$ws.Range("E35").Value = "  +3.79%  "

# Row 36 / 37: Bittensor <-> EthereumClassic swap position (with refreshed data)
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "18.52"
$ws.Range("E36").Value = "  +3.73%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "373.43"
$ws.Range("E37").Value = "  +15.96%  "

$ws.Range("D38").Value = "4.46"
$ws.Range("E38").Value = "  +9.41%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  +9.44%  "

$ws.Range("D42").Value = "39.06"
$ws.Range("E42").Value = "  +3.07%  "

$ws.Range("D43").Value = "145.33"
$ws.Range("E43").Value = "  +5.88%  "

$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  +4.79%  "

$ws.Range("D45").Value = "20.63"
$ws.Range("E45").Value = "  +7.48%  "

$ws.Range("E46").Value = "  +1.75%  "

# Row 47 / 48: Hedera <-> Mantle swap position (with refreshed data)
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.587"
$ws.Range("E47").Value = "  +3.98%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0518"
$ws.Range("E48").Value = "  +4.61%  "

$ws.Range("D49").Value = "18.00"
$ws.Range("E49").Value = "  +6.35%  "

$ws.Range("E50").Value = "  +3.31%  "

# Row 51: BabyDogeCoin replaced by dogwifhat
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  +11.39%  "

# ---------------------------------------------------------------------------
# Clear the temporary Text number-format back to the default "Normal" style
# so the cells end up with no explicit style index, matching the source
# workbook's original (styleless) string cells.
# ---------------------------------------------------------------------------
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
